$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---- Overview sheet: status columns (zh-cn / de-de) for both rows ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Range("I2").Value = "61c6aea6-efe2-48bd-beb1-b70d7a4b540d.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/181fd88429a293fe2f1d50782d86da3cd0414735/e2e/61c6aea6-efe2-48bd-beb1-b70d7a4b540d.md", "", "", "61c6aea6-efe2-48bd-beb1-b70d7a4b540d.md") | Out-Null
$wsZh.Range("J2").Value = "61c6aea6-efe2-48bd-beb1-b70d7a4b540d.1789728ae4b0d4115032646db52451364de813ad.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-16 10:24:15"

$wsZh.Range("I3").Value = "b35c1201-a9b3-422a-a4a1-5c715545fc4d.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/181fd88429a293fe2f1d50782d86da3cd0414735/e2e/b35c1201-a9b3-422a-a4a1-5c715545fc4d.md", "", "", "b35c1201-a9b3-422a-a4a1-5c715545fc4d.md") | Out-Null
$wsZh.Range("J3").Value = "b35c1201-a9b3-422a-a4a1-5c715545fc4d.0253029ecd842402ba0f09988d21a7a978965fa8.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-16 10:24:15"

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Range("I2").Value = "61c6aea6-efe2-48bd-beb1-b70d7a4b540d.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/181fd88429a293fe2f1d50782d86da3cd0414735/e2e/61c6aea6-efe2-48bd-beb1-b70d7a4b540d.md", "", "", "61c6aea6-efe2-48bd-beb1-b70d7a4b540d.md") | Out-Null
$wsDe.Range("J2").Value = "61c6aea6-efe2-48bd-beb1-b70d7a4b540d.1789728ae4b0d4115032646db52451364de813ad.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-16 10:24:22"

$wsDe.Range("I3").Value = "b35c1201-a9b3-422a-a4a1-5c715545fc4d.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/181fd88429a293fe2f1d50782d86da3cd0414735/e2e/b35c1201-a9b3-422a-a4a1-5c715545fc4d.md", "", "", "b35c1201-a9b3-422a-a4a1-5c715545fc4d.md") | Out-Null
$wsDe.Range("J3").Value = "b35c1201-a9b3-422a-a4a1-5c715545fc4d.0253029ecd842402ba0f09988d21a7a978965fa8.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-16 10:24:22"

# ---- Column width adjustments (to fit the new, longer text) ----
$wsOverview.Range("E1").ColumnWidth = 29.9777047293527
$wsOverview.Range("F1").ColumnWidth = 29.9777047293527

$wsZh.Range("C1").ColumnWidth = 29.9777047293527
$wsZh.Range("I1").ColumnWidth = 40
$wsZh.Range("J1").ColumnWidth = 40

$wsDe.Range("C1").ColumnWidth = 29.9777047293527
$wsDe.Range("I1").ColumnWidth = 40
$wsDe.Range("J1").ColumnWidth = 40
